$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.376.42"
$ws.Range("E2").Value = "  +1.89%  "

$ws.Range("D3").Value = "1.827.34"
$ws.Range("E3").Value = "  +1.08%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.02%  "

$ws.Range("E6").Value = "  -0.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4461"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.43%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3771"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.84%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07412"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.25%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8784"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.74%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.86"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.92%  "

$ws.Range("D12").Value = "1.828.62"
$ws.Range("E12").Value = "  +1.16%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.720"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.78%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.434"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.35%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.96"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.89%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07064"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.33%  "

$ws.Range("E17").Value = "  -0.04%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008820"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.18%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9999"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.04%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.41%  "

$ws.Range("D21").Value = "27.380.25"
$ws.Range("E21").Value = "  +1.85%  "

$ws.Range("E22").Value = "  +3.88%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.72%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.956"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.70%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.30"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.19%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.295"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.94%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.60"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.13%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.354"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.97%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.09"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.22%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08905"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.98%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.7961"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.87%  "

$ws.Range("E32").Value = "  +2.02%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.561"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.47%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.965"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.14%  "

$ws.Range("E36").Value = "  +1.28%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01979"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.59%  "

$ws.Range("E38").Value = "  +1.54%  "

$ws.Range("E39").Value = "  +3.65%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5343"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.42%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.872"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.30%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.338"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +18.61%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1705"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.99%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.703"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.08%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5080"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.87%  "

$ws.Range("E46").Value = "  +0.58%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "105.47"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.01%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.689"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.24%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.9994"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.03%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06396"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.93%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "66.12"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.53%  "
